$d = $word.ActiveDocument

# The "Notice" overview table is the first table in the document.
$t = $d.Tables.Item(1)

# --- Resize the table's four columns (widths given in twips in the
#     target OOXML; Word's COM object model reports/accepts Width in
#     points, so divide by 20). ---
$t.Columns.Item(1).Width = 2195 / 20
$t.Columns.Item(2).Width = 2445 / 20
$t.Columns.Item(3).Width = 528 / 20
$t.Columns.Item(4).Width = 2751 / 20

# --- Append two new notice rows at the bottom of the table. ---
$row1 = $t.Rows.Add()
$row1.Cells.Item(1).Range.Text = "no.elhub.flex.service_provider_product_suspension.product_type.not_qualified"
$row1.Cells.Item(2).Range.Text = "Inconsistency: suspending a SP on a product type that they are no longer qualified for"
$row1.Cells.Item(3).Range.Text = "PSO"
$row1.Cells.Item(4).Range.Text = "Delete the suspension because it is useless"

$row2 = $t.Rows.Add()
$row2.Cells.Item(1).Range.Text = "no.elhub.flex.service_provider_product_suspension.lingering"
$row2.Cells.Item(2).Range.Text = "Inactivity: nothing has happened on the suspension in 2 weeks"
$row2.Cells.Item(3).Range.Text = "PSO"
$row2.Cells.Item(4).Range.Text = "Suspension is a temporary procedure. Consider reinstating the SP or removing their qualification."

Write-Output ("Final row count=" + $t.Rows.Count)
